# Update "想去人数" (column F) figures across all sheets to match the
# refreshed data snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 12822
$ws1.Range("F3").Value  = 7172
$ws1.Range("F6").Value  = 451
$ws1.Range("F14").Value = 4
$ws1.Range("F19").Value = 370
$ws1.Range("F20").Value = 23
$ws1.Range("F21").Value = 279
$ws1.Range("F22").Value = 313
$ws1.Range("F24").Value = 168
$ws1.Range("F26").Value = 5242
$ws1.Range("F27").Value = 71
$ws1.Range("F28").Value = 1432
$ws1.Range("F29").Value = 313
$ws1.Range("F30").Value = 1375
$ws1.Range("F32").Value = 46
$ws1.Range("F33").Value = 1365
$ws1.Range("F38").Value = 3737

# Sheet: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 3744
$ws2.Range("F5").Value = 3744

# Sheet: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9283
$ws3.Range("F4").Value = 2023

# Sheet: 全部类型 (All Types) - combined view of all the above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 9283
$ws4.Range("F4").Value  = 2023
$ws4.Range("F5").Value  = 12822
$ws4.Range("F6").Value  = 7172
$ws4.Range("F8").Value  = 3744
$ws4.Range("F9").Value  = 451
$ws4.Range("F14").Value = 4
$ws4.Range("F19").Value = 370
$ws4.Range("F20").Value = 23
$ws4.Range("F21").Value = 279
$ws4.Range("F22").Value = 313
$ws4.Range("F27").Value = 168
$ws4.Range("F29").Value = 5242
$ws4.Range("F30").Value = 71
$ws4.Range("F31").Value = 1432
$ws4.Range("F34").Value = 313
$ws4.Range("F36").Value = 1375
$ws4.Range("F38").Value = 1365
$ws4.Range("F47").Value = 3737
